$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-13
$data = @(
    @{Row=2;  I=8; J=8},
    @{Row=3;  I=8; J=8},
    @{Row=4;  I=8; J=9},
    @{Row=5;  I=8; J=8},
    @{Row=6;  I=5; J=5},
    @{Row=7;  I=6; J=6},
    @{Row=8;  I=7; J=7},
    @{Row=9;  I=5; J=6},
    @{Row=10; I=9; J=9},
    @{Row=11; I=4; J=4},
    @{Row=12; I=8; J=8},
    @{Row=13; I=2; J=2}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}
